$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '67.846.17'
Set-TextValue $ws.Range('E2') '  -1.48%  '
Set-TextValue $ws.Range('D3') '3.798.90'
Set-TextValue $ws.Range('E3') '  -3.31%  '
Set-TextValue $ws.Range('E4') '  +0.25%  '
Set-TextValue $ws.Range('D5') '512.55'
Set-TextValue $ws.Range('E5') '  +4.87%  '
Set-TextValue $ws.Range('D6') '138.41'
Set-TextValue $ws.Range('E6') '  -5.56%  '
Set-TextValue $ws.Range('D7') '0.597'
Set-TextValue $ws.Range('E7') '  -4.20%  '
Set-TextValue $ws.Range('D8') '0.999'
Set-TextValue $ws.Range('E8') '  +0.14%  '
Set-TextValue $ws.Range('D9') '0.697'
Set-TextValue $ws.Range('E9') '  -5.61%  '
Set-TextValue $ws.Range('D10') '0.164'
Set-TextValue $ws.Range('E10') '  -7.43%  '
Set-TextValue $ws.Range('D11') '0.0000312'
Set-TextValue $ws.Range('E11') '  -10.88%  '
Set-TextValue $ws.Range('D12') '40.84'
Set-TextValue $ws.Range('E12') '  -4.83%  '
Set-TextValue $ws.Range('D13') '4.428.83'
Set-TextValue $ws.Range('E13') '  -2.84%  '
Set-TextValue $ws.Range('D14') '10.07'
Set-TextValue $ws.Range('E14') '  -4.01%  '
Set-TextValue $ws.Range('D15') '21.33'
Set-TextValue $ws.Range('E15') '  +6.86%  '
Set-TextValue $ws.Range('D16') '3.829.39'
Set-TextValue $ws.Range('E16') '  -2.55%  '
Set-TextValue $ws.Range('D17') '13.97'
Set-TextValue $ws.Range('E17') '  -2.21%  '
Set-TextValue $ws.Range('E18') '  -1.58%  '
Set-TextValue $ws.Range('D19') '1.18'
Set-TextValue $ws.Range('E19') '  +2.44%  '
Set-TextValue $ws.Range('D20') '67.989.19'
Set-TextValue $ws.Range('E20') '  -1.42%  '
Set-TextValue $ws.Range('D21') '410.46'
Set-TextValue $ws.Range('E21') '  -6.87%  '
Set-TextValue $ws.Range('D22') '3.34'
Set-TextValue $ws.Range('E22') '  -4.36%  '
Set-TextValue $ws.Range('D23') '13.80'
Set-TextValue $ws.Range('E23') '  -6.16%  '
Set-TextValue $ws.Range('D24') '85.58'
Set-TextValue $ws.Range('E24') '  -4.82%  '
Set-TextValue $ws.Range('B25') 'RenderToken'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D25') '11.51'
Set-TextValue $ws.Range('E25') '  -5.22%  '
Set-TextValue $ws.Range('B26') 'PancakeSwap'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D26') '3.87'
Set-TextValue $ws.Range('E26') '  +3.74%  '
Set-TextValue $ws.Range('D27') '10.32'
Set-TextValue $ws.Range('E27') '  -7.61%  '
Set-TextValue $ws.Range('D28') '34.81'
Set-TextValue $ws.Range('E28') '  -6.80%  '
Set-TextValue $ws.Range('D29') '677.82'
Set-TextValue $ws.Range('E29') '  -5.40%  '
Set-TextValue $ws.Range('D30') '12.94'
Set-TextValue $ws.Range('E30') '  -3.22%  '
Set-TextValue $ws.Range('D31') '0.123'
Set-TextValue $ws.Range('E31') '  -5.91%  '
Set-TextValue $ws.Range('D32') '2.78'
Set-TextValue $ws.Range('E32') '  -3.99%  '
Set-TextValue $ws.Range('D33') '64.42'
Set-TextValue $ws.Range('E33') '  +6.54%  '
Set-TextValue $ws.Range('B34') 'TheGraph'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D34') '0.442'
Set-TextValue $ws.Range('E34') '  -0.95%  '
Set-TextValue $ws.Range('D35') '5.94'
Set-TextValue $ws.Range('E35') '  -2.25%  '
Set-TextValue $ws.Range('B36') 'InjectiveProtocol'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D36') '39.04'
Set-TextValue $ws.Range('E36') '  -4.11%  '
Set-TextValue $ws.Range('B37') 'Dai'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D37') '0.999'
Set-TextValue $ws.Range('E37') '  +0.29%  '
Set-TextValue $ws.Range('B38') 'PEPE'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D38') '0.0₃0809'
Set-TextValue $ws.Range('E38') '  -9.47%  '
Set-TextValue $ws.Range('D39') '0.146'
Set-TextValue $ws.Range('E39') '  -1.47%  '
Set-TextValue $ws.Range('E40') '  -0.09%  '
Set-TextValue $ws.Range('D41') '3.23'
Set-TextValue $ws.Range('E41') '  +5.13%  '
Set-TextValue $ws.Range('D42') '0.0467'
Set-TextValue $ws.Range('E42') '  -4.66%  '
Set-TextValue $ws.Range('B43') 'Fetch.AI'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D43') '2.82'
Set-TextValue $ws.Range('E43') '  -4.68%  '
Set-TextValue $ws.Range('B44') 'WEMIXToken'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D44') '3.09'
Set-TextValue $ws.Range('E44') '  +3.50%  '
Set-TextValue $ws.Range('D45') '3.33'
Set-TextValue $ws.Range('E45') '  -2.95%  '
Set-TextValue $ws.Range('D46') '0.137'
Set-TextValue $ws.Range('E46') '  -3.98%  '
Set-TextValue $ws.Range('D47') '2.89'
Set-TextValue $ws.Range('E47') '  -1.03%  '
Set-TextValue $ws.Range('D48') '143.04'
Set-TextValue $ws.Range('E48') '  +0.68%  '
Set-TextValue $ws.Range('D49') '2.671.87'
Set-TextValue $ws.Range('E49') '  +9.18%  '
Set-TextValue $ws.Range('B50') 'LidoDAOToken'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D50') '3.20'
Set-TextValue $ws.Range('E50') '  -5.38%  '
Set-TextValue $ws.Range('B51') 'ARBITRUM'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D51') '2.01'
Set-TextValue $ws.Range('E51') '  -3.03%  '
